$wb = $excel.ActiveWorkbook

# --- Sheet2: keep same text values, just re-pointed at the (now lower) shared-string
#     indices once the stale "Test outlet"/"St2235" strings are pruned below. The
#     leading apostrophe forces these numeric-looking strings to stay text (t="s").
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("A2").Value = "'654645"
$ws2.Range("B2").Value = "'6466546"

# --- StatementCollect: new header column + two data rows
#     (write order matches the shared-string append order: MZ-0004, STA00006,
#     MZ-0005, STA00007, then CollectAmt last)
$ws3 = $wb.Worksheets.Item("StatementCollect")

$ws3.Range("A2").Value = "MZ-0004"
$ws3.Range("B2").Value = "STA00006"
$ws3.Range("C2").Value = 11

$ws3.Range("A3").Value = "MZ-0005"
$ws3.Range("B3").Value = "STA00007"
$ws3.Range("C3").Value = 2

$ws3.Range("C1").Value = "CollectAmt"

# StatementCollect becomes the active sheet/tab, with G4 selected
$ws3.Activate() | Out-Null
$ws3.Range("G4").Select() | Out-Null
